# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.881.11'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.414.09'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '551.60'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.01'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.591'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.97%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.70'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.71%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.15%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.60'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.842.64'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '59.804.97'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.378.36'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.32'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '328.42'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.45%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.67%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.43'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.172'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.04%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.00%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.11%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.38%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.11'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.42%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.63%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.43%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.00%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.21'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.23%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '313.61'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.28%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.49%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '138.69'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0969'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0516'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.51'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.02%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0223'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '17.64'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.28%  '
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.10'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.57'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.67%  '
$ws.Range("B51").Value = 'ZEEBU'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.67'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.11%  '
